$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32 (@@ -2209,25 +2209,22 @@)
$ws.Range("H32").Value = 4300
$ws.Range("I32").Value = 4300
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4300
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3974
$ws.Range("N32").ClearContents()

# Row 76 (@@ -4377,20 +4374,23 @@)
$ws.Range("H76").Value = 3999.6667
$ws.Range("I76").Value = 3999
$ws.Range("K76").Value = 3999
$ws.Range("M76").Value = -3684

# Row 79 (@@ -4530,20 +4530,23 @@)
$ws.Range("H79").Value = 3999.6667
$ws.Range("I79").Value = 3999
$ws.Range("K79").Value = 3999
$ws.Range("M79").Value = -2907

# Row 80 (@@ -4579,25 +4582,25 @@)
$ws.Range("H80").Value = 908.8889
$ws.Range("J80").Value = 1186.75
$ws.Range("L80").Value = 3560.25
$ws.Range("N80").Value = -5556.25

# Row 83 (@@ -4726,25 +4729,25 @@)
$ws.Range("H83").Value = 908.8889
$ws.Range("J83").Value = 1186.75
$ws.Range("L83").Value = 10680.75
$ws.Range("N83").Value = -20664.75

# Row 103 (@@ -5715,25 +5718,25 @@)
$ws.Range("H103").Value = 1931.5555
$ws.Range("I103").Value = 2838.3333
$ws.Range("J103").Value = 1478.1666
$ws.Range("K103").Value = 8514.999899999999
$ws.Range("L103").Value = 4434.4998
$ws.Range("M103").Value = -7928.999899999999
$ws.Range("N103").Value = -5606.4998

# Row 135 (@@ -7280,25 +7283,25 @@)
$ws.Range("H135").Value = 692.2308
$ws.Range("I135").Value = 834.2222
$ws.Range("J135").Value = 372.75
$ws.Range("K135").Value = 7507.999800000001
$ws.Range("L135").Value = 3354.75
$ws.Range("M135").Value = -4972.999800000001
$ws.Range("N135").Value = -8424.75

# Row 141 (@@ -7580,22 +7583,22 @@)
$ws.Range("H141").Value = 3365
$ws.Range("I141").Value = 997.5
$ws.Range("K141").Value = 2992.5
$ws.Range("M141").Value = 2187.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (@@ -7726,25 +7729,25 @@)
$ws.Range("H2").Value = 935.7
$ws.Range("I2").Value = 844.25
$ws.Range("J2").Value = 1072.875
$ws.Range("K2").Value = 844.25
$ws.Range("L2").Value = 1072.875
$ws.Range("M2").Value = -731.25
$ws.Range("N2").Value = -1298.875

# Row 116 (@@ -13279,25 +13282,25 @@)
$ws.Range("H116").Value = 935.7
$ws.Range("I116").Value = 844.25
$ws.Range("J116").Value = 1072.875
$ws.Range("K116").Value = 844.25
$ws.Range("L116").Value = 1072.875
$ws.Range("M116").Value = 1449.75
$ws.Range("N116").Value = -5660.875

# Row 122 (@@ -13573,22 +13576,22 @@)
$ws.Range("H122").Value = 1657.8
$ws.Range("I122").Value = 1296.0435
$ws.Range("K122").Value = 3888.1305
$ws.Range("M122").Value = -1438.1305

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (@@ -14675,25 +14678,25 @@)
$ws.Range("H3").Value = 935.7
$ws.Range("I3").Value = 844.25
$ws.Range("J3").Value = 1072.875
$ws.Range("K3").Value = 844.25
$ws.Range("L3").Value = 1072.875
$ws.Range("M3").Value = -730.25
$ws.Range("N3").Value = -1300.875

# Row 95 (@@ -19117,22 +19120,22 @@)
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492

# Row 107 (@@ -19696,25 +19699,25 @@)
$ws.Range("H107").Value = 4078.3333
$ws.Range("I107").Value = 3464.1428
$ws.Range("J107").Value = 4615.75
$ws.Range("K107").Value = 3464.1428
$ws.Range("L107").Value = 4615.75
$ws.Range("M107").Value = -1544.1428
$ws.Range("N107").Value = -8455.75

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (@@ -21699,22 +21702,22 @@)
$ws.Range("H7").Value = 3554.9656
$ws.Range("I7").Value = 5376.9473
$ws.Range("K7").Value = 5376.9473
$ws.Range("M7").Value = -5263.9473

# Row 99 (@@ -26189,22 +26192,22 @@)
$ws.Range("H99").Value = 2195
$ws.Range("I99").Value = 2195
$ws.Range("K99").Value = 2195
$ws.Range("M99").Value = -697

# Row 122 (@@ -27307,22 +27310,22 @@)
$ws.Range("H122").Value = 4998
$ws.Range("I122").Value = 4997
$ws.Range("K122").Value = 14991
$ws.Range("M122").Value = -12541

# Row 126 (@@ -27500,22 +27503,22 @@)
$ws.Range("H126").Value = 2195
$ws.Range("I126").Value = 2195
$ws.Range("K126").Value = 6585
$ws.Range("M126").Value = -4115

$ws = $wb.Worksheets.Item("CUL")
# Row 23 (@@ -29437,25 +29440,25 @@)
$ws.Range("H23").Value = 141.5
$ws.Range("J23").Value = 166.66667
$ws.Range("L23").Value = 500.00001
$ws.Range("N23").Value = -970.00001

# Row 86 (@@ -32590,25 +32593,25 @@)
$ws.Range("H86").Value = 539.5333000000001
$ws.Range("J86").Value = 564.3
$ws.Range("L86").Value = 1692.9
$ws.Range("N86").Value = -4064.9

# Row 89 (@@ -32740,25 +32743,25 @@)
$ws.Range("H89").Value = 539.5333000000001
$ws.Range("J89").Value = 564.3
$ws.Range("L89").Value = 5078.7
$ws.Range("N89").Value = -16934.7

# Row 113 (@@ -33925,25 +33928,25 @@)
$ws.Range("H113").Value = 560.6
$ws.Range("J113").Value = 640.2
$ws.Range("L113").Value = 1920.6
$ws.Range("N113").Value = -6260.6

$ws = $wb.Worksheets.Item("GSM")
# Row 39 (@@ -37295,22 +37298,19 @@)
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

# Row 49 (@@ -37779,22 +37779,19 @@)
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

# Row 102 (@@ -40319,22 +40316,22 @@)
$ws.Range("H102").Value = 1492.2727
$ws.Range("I102").Value = 1492.2727
$ws.Range("K102").Value = 1492.2727
$ws.Range("M102").Value = 129.7273

# Row 122 (@@ -41272,22 +41269,22 @@)
$ws.Range("H122").Value = 1813.1666
$ws.Range("I122").Value = 1813.1666
$ws.Range("K122").Value = 5439.4998
$ws.Range("M122").Value = -2989.4998

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (@@ -43296,25 +43293,25 @@)
$ws.Range("H22").Value = 865.63635
$ws.Range("I22").Value = 835.7778
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 835.7778
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -540.7778
$ws.Range("N22").Value = -1590

# Row 27 (@@ -43541,25 +43538,25 @@)
$ws.Range("H27").Value = 865.63635
$ws.Range("I27").Value = 835.7778
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 835.7778
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -728.7778
$ws.Range("N27").Value = -1214

# Row 40 (@@ -44175,22 +44172,19 @@)
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

# Row 46 (@@ -44466,25 +44460,25 @@)
$ws.Range("H46").Value = 5399.231
$ws.Range("J46").Value = 6885.7144
$ws.Range("L46").Value = 6885.7144
$ws.Range("N46").Value = -7261.7144

# Row 55 (@@ -44889,25 +44883,25 @@)
$ws.Range("H55").Value = 1338.2273
$ws.Range("I55").Value = 928.8
$ws.Range("J55").Value = 1679.4166
$ws.Range("K55").Value = 928.8
$ws.Range("L55").Value = 1679.4166
$ws.Range("M55").Value = -755.8
$ws.Range("N55").Value = -2025.4166

# Row 68 (@@ -45520,25 +45514,25 @@)
$ws.Range("H68").Value = 8665.799999999999
$ws.Range("I68").Value = 8533
$ws.Range("J68").Value = 8798.6
$ws.Range("K68").Value = 8533
$ws.Range("L68").Value = 8798.6
$ws.Range("M68").Value = -7784
$ws.Range("N68").Value = -10296.6

# Row 71 (@@ -45664,25 +45658,25 @@)
$ws.Range("H71").Value = 8665.799999999999
$ws.Range("I71").Value = 8533
$ws.Range("J71").Value = 8798.6
$ws.Range("K71").Value = 42665
$ws.Range("L71").Value = 43993
$ws.Range("M71").Value = -38921
$ws.Range("N71").Value = -51481

# Row 82 (@@ -46182,25 +46176,25 @@)
$ws.Range("H82").Value = 6927.857
$ws.Range("J82").Value = 8000
$ws.Range("L82").Value = 8000
$ws.Range("N82").Value = -8722

# Row 85 (@@ -46326,25 +46320,25 @@)
$ws.Range("H85").Value = 6927.857
$ws.Range("J85").Value = 8000
$ws.Range("L85").Value = 8000
$ws.Range("N85").Value = -10496

# Row 93 (@@ -46706,22 +46700,22 @@)
$ws.Range("H93").Value = 998
$ws.Range("I93").Value = 997
$ws.Range("K93").Value = 997
$ws.Range("M93").Value = 251

# Row 97 (@@ -46899,19 +46893,22 @@)
$ws.Range("H97").Value = 20000
$ws.Range("J97").Value = 20000
$ws.Range("L97").Value = 20000
$ws.Range("N97").Value = -21982

# Row 122 (@@ -48088,19 +48085,22 @@)
$ws.Range("H122").Value = 2004
$ws.Range("I122").Value = 2004
$ws.Range("K122").Value = 6012
$ws.Range("M122").Value = -3562

# Row 132 (@@ -48560,22 +48560,22 @@)
$ws.Range("H132").Value = 3994.5
$ws.Range("I132").Value = 3994.5
$ws.Range("K132").Value = 11983.5
$ws.Range("M132").Value = -9453.5

$ws = $wb.Worksheets.Item("WVR")
# Row 41 (@@ -51007,25 +51007,25 @@)
$ws.Range("H41").Value = 34119.2
$ws.Range("J41").Value = 42683.5
$ws.Range("L41").Value = 42683.5
$ws.Range("N41").Value = -43463.5

# Row 45 (@@ -51197,25 +51197,25 @@)
$ws.Range("H45").Value = 25541.666
$ws.Range("J45").Value = 29313
$ws.Range("L45").Value = 29313
$ws.Range("N45").Value = -30295
